$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so numeric-looking
# strings (e.g. "0.593") are not auto-converted to numbers by Excel,
# matching the source data which stores these as plain text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.451.14'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '2.716.07'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '558.39'
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("D6").Value = '157.47'
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("E9").Value = '  -3.08%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '5.60'
$ws.Range("E11").Value = '  -5.22%  '
$ws.Range("D12").Value = '0.371'
$ws.Range("E12").Value = '  -4.58%  '
$ws.Range("D13").Value = '3.196.52'
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").Value = '26.51'
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = '63.317.95'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("E16").Value = '  -3.73%  '
$ws.Range("D17").Value = '2.719.90'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '12.13'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").Value = '4.64'
$ws.Range("E19").Value = '  -4.85%  '
$ws.Range("D20").Value = '349.53'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").Value = '6.41'
$ws.Range("E21").Value = '  -4.89%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '0.513'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("D24").Value = '64.33'
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = '8.16'
$ws.Range("E27").Value = '  -5.95%  '
$ws.Range("D28").Value = '0.0₃0878'
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '1.95'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  +7.17%  '
$ws.Range("E31").Value = '  -2.92%  '
$ws.Range("D32").Value = '164.51'
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '19.84'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.47'
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").Value = '4.80'
$ws.Range("E36").Value = '  -3.32%  '
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  -0.76%  '
$ws.Range("D38").Value = '346.66'
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("D39").Value = '0.955'
$ws.Range("E39").Value = '  -5.25%  '
$ws.Range("D40").Value = '6.02'
$ws.Range("E40").Value = '  -4.79%  '
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("D42").Value = '38.25'
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").Value = '21.25'
$ws.Range("E43").Value = '  -3.55%  '
$ws.Range("D44").Value = '20.65'
$ws.Range("E44").Value = '  -4.67%  '
$ws.Range("D45").Value = '0.0569'
$ws.Range("E45").Value = '  -3.80%  '
$ws.Range("D46").Value = '0.624'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '131.83'
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '11.07'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '0.0982'
$ws.Range("E50").Value = '  -3.80%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '0.0244'
$ws.Range("E51").Value = '  -4.52%  '
